$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute(
    "Content Work for VA.gov Brand Consolidation: Vets.gov",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Content Work for VA.gov Brand Consolidation",
    2
)
